$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Update usernames from placeholder values to real Ukrainian names
$ws.Range("A2").Value = "Іваненко"
$ws.Range("A3").Value = "Петренко"
$ws.Range("A4").Value = "Заічко"

# Update the active cell selection
$ws.Range("A8").Select()
